$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.6275843652054141
$ws.Range("J2").Value = 0.6275843652054141
$ws.Range("M2").Value = 0.9949089999999999
$ws.Range("N2").Value = 2.984727
$ws.Range("O2").Value = 0.1476822527339178
$ws.Range("P2").Value = 0.1476822527339178
$ws.Range("Q2").Value = 0.014808225556
$ws.Range("R2").Value = 0.133274030004
$ws.Range("S2").Value = 0.09268307283412133
$ws.Range("T2").Value = 0.09268307283412133
$ws.Range("I3").Value = 0.6275843652054141
$ws.Range("J3").Value = 0.6275843652054141
$ws.Range("O3").Value = 0.2453919293791607
$ws.Range("P3").Value = 0.2453919293791607
$ws.Range("S3").Value = 0.1540041382259524
$ws.Range("T3").Value = 0.1540041382259524
$ws.Range("I4").Value = 0.6275843652054141
$ws.Range("J4").Value = 0.6275843652054141
$ws.Range("M4").Value = 1.748891
$ws.Range("N4").Value = 5.246673
$ws.Range("O4").Value = 0.2596017954064887
$ws.Range("P4").Value = 0.2596017954064887
$ws.Range("Q4").Value = 0.026030493644
$ws.Range("R4").Value = 0.234274442796
$ws.Range("S4").Value = 0.162922027976367
$ws.Range("T4").Value = 0.162922027976367
$ws.Range("I5").Value = 0.6275843652054141
$ws.Range("J5").Value = 0.6275843652054141
$ws.Range("M5").Value = 0.7268083333333334
$ws.Range("N5").Value = 2.180425
$ws.Range("O5").Value = 0.1078859392893731
$ws.Range("P5").Value = 0.1078859392893731
$ws.Range("Q5").Value = 0.01081781523333333
$ws.Range("R5").Value = 0.0973603371
$ws.Range("S5").Value = 0.06770752872351106
$ws.Range("T5").Value = 0.06770752872351106
$ws.Range("I6").Value = 0.6275843652054141
$ws.Range("J6").Value = 0.6275843652054141
$ws.Range("M6").Value = 1.613051666666667
$ws.Range("N6").Value = 4.839155
$ws.Range("O6").Value = 0.2394380831910597
$ws.Range("P6").Value = 0.2394380831910597
$ws.Range("Q6").Value = 0.02400866100666666
$ws.Range("R6").Value = 0.21607794906
$ws.Range("S6").Value = 0.1502675974454623
$ws.Range("T6").Value = 0.1502675974454623
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.008832333333333333
$ws.Range("H7").Value = 0.026497
$ws.Range("I7").Value = 0.372415634794586
$ws.Range("J7").Value = 0.3724156347945861
$ws.Range("M7").Value = 0.9949089999999999
$ws.Range("N7").Value = 2.984727
$ws.Range("O7").Value = 0.1476822527339178
$ws.Range("P7").Value = 0.1476822527339178
$ws.Range("Q7").Value = 0.008787367924333333
$ws.Range("R7").Value = 0.079086311319
$ws.Range("S7").Value = 0.05499917989979648
$ws.Range("T7").Value = 0.05499917989979649
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.008832333333333333
$ws.Range("H8").Value = 0.026497
$ws.Range("I8").Value = 0.372415634794586
$ws.Range("J8").Value = 0.3724156347945861
$ws.Range("O8").Value = 0.2453919293791607
$ws.Range("P8").Value = 0.2453919293791607
$ws.Range("Q8").Value = 0.01460127489388889
$ws.Range("R8").Value = 0.131411474045
$ws.Range("S8").Value = 0.09138779115320836
$ws.Range("T8").Value = 0.09138779115320837
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.008832333333333333
$ws.Range("H9").Value = 0.026497
$ws.Range("I9").Value = 0.372415634794586
$ws.Range("J9").Value = 0.3724156347945861
$ws.Range("M9").Value = 1.748891
$ws.Range("N9").Value = 5.246673
$ws.Range("O9").Value = 0.2596017954064887
$ws.Range("P9").Value = 0.2596017954064887
$ws.Range("Q9").Value = 0.01544678827566667
$ws.Range("R9").Value = 0.139021094481
$ws.Range("S9").Value = 0.09667976743012173
$ws.Range("T9").Value = 0.09667976743012174
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.008832333333333333
$ws.Range("H10").Value = 0.026497
$ws.Range("I10").Value = 0.372415634794586
$ws.Range("J10").Value = 0.3724156347945861
$ws.Range("M10").Value = 0.7268083333333334
$ws.Range("N10").Value = 2.180425
$ws.Range("O10").Value = 0.1078859392893731
$ws.Range("P10").Value = 0.1078859392893731
$ws.Range("Q10").Value = 0.006419413469444445
$ws.Range("R10").Value = 0.057774721225
$ws.Range("S10").Value = 0.04017841056586206
$ws.Range("T10").Value = 0.04017841056586206
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.008832333333333333
$ws.Range("H11").Value = 0.026497
$ws.Range("I11").Value = 0.372415634794586
$ws.Range("J11").Value = 0.3724156347945861
$ws.Range("M11").Value = 1.613051666666667
$ws.Range("N11").Value = 4.839155
$ws.Range("O11").Value = 0.2394380831910597
$ws.Range("P11").Value = 0.2394380831910597
$ws.Range("Q11").Value = 0.01424701000388889
$ws.Range("R11").Value = 0.128223090035
$ws.Range("S11").Value = 0.08917048574559738
$ws.Range("T11").Value = 0.0891704857455974
